$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.287.28'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '1.864.24'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''236.45'
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '''0.4708'
$ws.Range("E7").Value = '  +1.01%  '
$ws.Range("D8").Value = '''0.2902'
$ws.Range("E8").Value = '  +2.21%  '
$ws.Range("D9").Value = '''0.06541'
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").Value = '''21.90'
$ws.Range("E10").Value = '  +3.12%  '
$ws.Range("D11").Value = '''0.07946'
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").Value = '''97.86'
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = '1.877.43'
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("D14").Value = '''5.151'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '''0.6808'
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("D16").Value = '''263.16'
$ws.Range("E16").Value = '  -6.43%  '
$ws.Range("D17").Value = '30.265.07'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '''13.73'
$ws.Range("E18").Value = '  +8.31%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '''0.000007471'
$ws.Range("E20").Value = '  +2.50%  '
$ws.Range("D21").Value = '2.112.56'
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '''5.277'
$ws.Range("E23").Value = '  -4.48%  '
$ws.Range("D24").Value = '''6.182'
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").Value = '''167.34'
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("D26").Value = '''9.192'
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("D27").Value = '''18.91'
$ws.Range("E27").Value = '  -1.85%  '
$ws.Range("D28").Value = '''1.951'
$ws.Range("E28").Value = '  +1.17%  '
$ws.Range("D29").Value = '''1.395'
$ws.Range("E29").Value = '  +1.59%  '
$ws.Range("D30").Value = '''0.09883'
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D31").Value = '''4.355'
$ws.Range("E31").Value = '  -1.56%  '
$ws.Range("D32").Value = '''1.472'
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("D34").Value = '''0.04717'
$ws.Range("E34").Value = '  +0.61%  '
$ws.Range("D35").Value = '''1.131'
$ws.Range("E35").Value = '  +0.79%  '
$ws.Range("D36").Value = '''0.7014'
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("D37").Value = '''2.709'
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("D38").Value = '''0.01882'
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("D39").Value = '''2.621'
$ws.Range("E39").Value = '  +3.23%  '
$ws.Range("D40").Value = '''6.326'
$ws.Range("E40").Value = '  +0.70%  '
$ws.Range("D41").Value = '''73.90'
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("D42").Value = '''1.944'
$ws.Range("E42").Value = '  -0.34%  '
$ws.Range("D43").Value = '''0.8430'
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.4165'
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '''1.000'
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '''103.34'
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("D48").Value = '''951.20'
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("D49").Value = '''9.228'
$ws.Range("E49").Value = '  +0.96%  '
$ws.Range("D50").Value = '''34.18'
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("D51").Value = '''0.05664'
$ws.Range("E51").Value = '  +0.63%  '

# Reset style on cells forced to text via leading apostrophe, to avoid
# leaving a stray NumberFormat/quotePrefix flag on the cell.
$ws.Range("D5").Style = 'Normal'
$ws.Range("D6").Style = 'Normal'
$ws.Range("D7").Style = 'Normal'
$ws.Range("D8").Style = 'Normal'
$ws.Range("D9").Style = 'Normal'
$ws.Range("D10").Style = 'Normal'
$ws.Range("D11").Style = 'Normal'
$ws.Range("D12").Style = 'Normal'
$ws.Range("D14").Style = 'Normal'
$ws.Range("D15").Style = 'Normal'
$ws.Range("D16").Style = 'Normal'
$ws.Range("D18").Style = 'Normal'
$ws.Range("D20").Style = 'Normal'
$ws.Range("D23").Style = 'Normal'
$ws.Range("D24").Style = 'Normal'
$ws.Range("D25").Style = 'Normal'
$ws.Range("D26").Style = 'Normal'
$ws.Range("D27").Style = 'Normal'
$ws.Range("D28").Style = 'Normal'
$ws.Range("D29").Style = 'Normal'
$ws.Range("D30").Style = 'Normal'
$ws.Range("D31").Style = 'Normal'
$ws.Range("D32").Style = 'Normal'
$ws.Range("D34").Style = 'Normal'
$ws.Range("D35").Style = 'Normal'
$ws.Range("D36").Style = 'Normal'
$ws.Range("D37").Style = 'Normal'
$ws.Range("D38").Style = 'Normal'
$ws.Range("D39").Style = 'Normal'
$ws.Range("D40").Style = 'Normal'
$ws.Range("D41").Style = 'Normal'
$ws.Range("D42").Style = 'Normal'
$ws.Range("D43").Style = 'Normal'
$ws.Range("D44").Style = 'Normal'
$ws.Range("D45").Style = 'Normal'
$ws.Range("D46").Style = 'Normal'
$ws.Range("D48").Style = 'Normal'
$ws.Range("D49").Style = 'Normal'
$ws.Range("D50").Style = 'Normal'
$ws.Range("D51").Style = 'Normal'
